$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.324.03"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "1.653.24"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "1.889.75"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "1.651.79"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.569"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "27.340.56"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.59%  "
$ws.Range("D19").Value = "0.0₃0728"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "1.456.17"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.10"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.38"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.907"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.570"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.50%  "
$ws.Range("D45").Value = "1.797.02"
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.785"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("E49").Value = "  -3.35%  "
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.30%  "
